# NIT-9001494600.xlsx update: refresh "Estado de Cuenta" employee detail
# table with the latest database extract (Leonardo Tovio 1901-1907,
# Shirley Hernandez 1902-1907, Santiago Padilla 2507-2508) and drop the
# employees/periods that are no longer in arrears (Adolfo Garces, Thamara
# De Avila, Dianis Rodriguez, and Santiago's old 2505/2506/2507 periods).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the old detail rows that no longer belong in the statement
#    (old rows 31-38: Thamara De Avila x2, Dianis Rodriguez x3, Santiago
#    Padilla's previous 2505/2506/2507 x3). Deleting these entire rows
#    shifts the blank spacer + signature block (old rows 39-44) up by 8
#    rows automatically (-> new rows 31-36).
$ws.Range("B31:B38").EntireRow.Delete()

# 2) Rewrite the remaining detail rows (16-30) with the refreshed data.
#    Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
#    E=Periodo Mora, F=Valor Mora, G=Salario Basico.
$rows = @(
  @{ Row = 16; B = "CC"; C = "1051824735"; D = "LEONARDO RAFAEL TOVIO OSORIO";        E = "1901"; F = 58000;  G = 2117100 },
  @{ Row = 17; B = "CC"; C = "1051824735"; D = "LEONARDO RAFAEL TOVIO OSORIO";        E = "1902"; F = 60000;  G = 2117100 },
  @{ Row = 18; B = "CC"; C = "1143339361"; D = "SHIRLEY JISETH HERNANDEZ CABEZA";     E = "1902"; F = 32000;  G = 4760000 },
  @{ Row = 19; B = "CC"; C = "1051824735"; D = "LEONARDO RAFAEL TOVIO OSORIO";        E = "1903"; F = 60000;  G = 2117100 },
  @{ Row = 20; B = "CC"; C = "1143339361"; D = "SHIRLEY JISETH HERNANDEZ CABEZA";     E = "1903"; F = 160000; G = 4760000 },
  @{ Row = 21; B = "CC"; C = "1051824735"; D = "LEONARDO RAFAEL TOVIO OSORIO";        E = "1904"; F = 60000;  G = 2117100 },
  @{ Row = 22; B = "CC"; C = "1143339361"; D = "SHIRLEY JISETH HERNANDEZ CABEZA";     E = "1904"; F = 160000; G = 4760000 },
  @{ Row = 23; B = "CC"; C = "1051824735"; D = "LEONARDO RAFAEL TOVIO OSORIO";        E = "1905"; F = 60000;  G = 2117100 },
  @{ Row = 24; B = "CC"; C = "1143339361"; D = "SHIRLEY JISETH HERNANDEZ CABEZA";     E = "1905"; F = 160000; G = 4760000 },
  @{ Row = 25; B = "CC"; C = "1051824735"; D = "LEONARDO RAFAEL TOVIO OSORIO";        E = "1906"; F = 60000;  G = 2117100 },
  @{ Row = 26; B = "CC"; C = "1143339361"; D = "SHIRLEY JISETH HERNANDEZ CABEZA";     E = "1906"; F = 160000; G = 4760000 },
  @{ Row = 27; B = "CC"; C = "1051824735"; D = "LEONARDO RAFAEL TOVIO OSORIO";        E = "1907"; F = 60000;  G = 2117100 },
  @{ Row = 28; B = "CC"; C = "1143339361"; D = "SHIRLEY JISETH HERNANDEZ CABEZA";     E = "1907"; F = 160000; G = 4760000 },
  @{ Row = 29; B = "CC"; C = "1007315093"; D = "SANTIAGO DE JESUS PADILLA MARMOLEJO"; E = "2507"; F = 28470;  G = 711750 },
  @{ Row = 30; B = "CC"; C = "1007315093"; D = "SANTIAGO DE JESUS PADILLA MARMOLEJO"; E = "2508"; F = 28470;  G = 711750 }
)

foreach ($row in $rows) {
  $r = $row.Row
  $ws.Cells.Item($r, 2).Value = $row.B
  $ws.Cells.Item($r, 3).Value = $row.C
  $ws.Cells.Item($r, 4).Value = $row.D
  $ws.Cells.Item($r, 5).Value = $row.E
  $ws.Cells.Item($r, 6).Value = $row.F
  $ws.Cells.Item($r, 7).Value = $row.G
}

# 3) Refresh the summary figures at the top of the statement.
$ws.Range("E11").Value = 1306940   # VALOR MORA (sum of the 15 remaining rows)
$ws.Range("C13").Value = 3         # Cant. Trabajadores
$ws.Range("F13").Value = 9         # Cant. Periodos
